# daily auto push: 2026-02-18 03:16 UTC
# Insert a new data row (2026/02/18, 水, 7, 201) right before the existing
# row that currently holds "2026/12/29" (row 827), shifting every row
# below it down by one. The sheet's used range grows from A1:D868 to
# A1:D869.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 827

# Push row 827..868 down to 828..869 by inserting a blank row above 827.
$ws.Rows.Item($insertRow).Insert()

# Populate the freshly inserted row with the new record. The date column
# stores plain text (e.g. "2026/02/18"), not a real date value, so force a
# text format while assigning it and then drop back to the sheet's default
# (unstyled) look, matching every other data row.
$ws.Cells.Item($insertRow, 1).NumberFormat = "@"
$ws.Cells.Item($insertRow, 1).Value = "2026/02/18"
$ws.Cells.Item($insertRow, 1).Style = "Normal"

$ws.Cells.Item($insertRow, 2).Value = "水"
$ws.Cells.Item($insertRow, 3).Value = 7
$ws.Cells.Item($insertRow, 4).Value = 201
